$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $new, 2) | Out-Null
}

# First occurrence block (Original Output)
Replace-Text "arr (prgram heap) location: 0000000CF82FF6C8" "arr (prgram heap- lifetime: program end) location: 0000000CF82FF6C8"
Replace-Text "x (static) location: 00007FF6CF07E200" "x (static- lifetime: program end) location: 00007FF6CF07E200"
Replace-Text "y (code) location: 00007FF6CF07E228" "y (code- lifetime: function end) location: 00007FF6CF07E228"
Replace-Text "hello function (code) 00007FF6CF07133E" "hello function (code- lifetime: program end) 00007FF6CF07133E"
Replace-Text "next function (code) 00007FF6CF071154" "next function (code- lifetime: program end) 00007FF6CF071154"
Replace-Text "inc (call stack) location: 0000000CF82FF5A4" "inc (call stack- lifetime: function end) location: 0000000CF82FF5A4"
Replace-Text "j (static) location: 00007FF6CF07E000" "j (static- lifetime: program end) location: 00007FF6CF07E000"
Replace-Text "i location: 0000000CF82FF6A0" "i (call stack- lifetime: function end) location: 0000000CF82FF6A0"
